$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "SingleBureauSingleLineVisa"
$ws1.Range("C2").Value = "4716428411725021"
